{"js": "// Insert a new \"2 F\u00c9VRIER\" row at the very top of the first table (the\n// \"JOURNAL DE BORD\" log table), mirroring the existing row layout: a date\n// cell on the left and a bulleted (\"Paragraphedeliste\"/numId 2) entry on\n// the right.\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst dateText = \"2 F\u00c9VRIER\";\nconst entryText =\n  \"Cr\u00e9ation des formulaires d\\u2019entrevue et des projets de stages et \" +\n  \"changement du formulaire d\\u2019\u00e9valuation des coordonnateurs en html. \" +\n  \"Tous ces formulaires sont maintenant en double soient en lecture et en \" +\n  \"\u00e9criture.\";\n\n// addRows(\"Start\", \u2026) inserts before the current first row and copies that\n// row's paragraph/list formatting (bullet list on the right-hand cell),\n// matching the rest of the table.\ntable.addRows(\"Start\", 1, [[dateText, entryText]]);\nawait context.sync();\n", "ps1": "# Insert a new \"2 F\u00c9VRIER\" row at the very top of the first table (the\n# \"JOURNAL DE BORD\" log table), mirroring the existing row layout: a date\n# cell on the left and a bulleted entry on the right.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$firstRow = $t.Rows.Item(1)\n# Rows.Add(BeforeRow) inserts a new row immediately above $firstRow, copying\n# that row's paragraph/list formatting (bullet list on the right-hand cell),\n# matching the rest of the table.\n$newRow = $t.Rows.Add($firstRow)\n\n$newRow.Cells.Item(1).Range.Text = \"2 F\u00c9VRIER\"\n$newRow.Cells.Item(2).Range.Text = \"Cr\u00e9ation des formulaires d\u2019entrevue et des projets de stages et changement du formulaire d\u2019\u00e9valuation des coordonnateurs en html. Tous ces formulaires sont maintenant en double soient en lecture et en \u00e9criture.\"\n"}
